$wb = $excel.ActiveWorkbook

# --- Rename the original (only) sheet to "E3" ---
$orig = $wb.Worksheets.Item(1)
$orig.Name = "E3"

# --- Duplicate "E3" and place the copy before it; rename the copy "E3 blank" ---
$orig.Copy($orig)
$wb.Worksheets.Item(1).Name = "E3 blank"

# Re-fetch both sheets by name so we always hold the correct, live objects
# (the copy operation can otherwise leave stale references pointing at the
# wrong sheet).
$e3 = $wb.Worksheets.Item("E3")
$blank = $wb.Worksheets.Item("E3 blank")

# --- On the "E3 blank" sheet, clear out the student-entry / computed cells ---
# (keeps labels/headers, but empties the volume inputs and every derived value)
$blank.Range("B2:C2").ClearContents()
$blank.Range("B3:C4").ClearContents()
$blank.Range("B5").ClearContents()
$blank.Range("B8:C9").ClearContents()

# --- Highlight the constraint-formula cells in red on both sheets ---
$blank.Range("B8:C9").Interior.Color = 255
$e3.Range("B8:C9").Interior.Color = 255

# --- Make "E3 blank" the active/selected tab, matching the authored file ---
# (selecting a range also activates its sheet, so set "E3" first and finish
# on "E3 blank" so it is the one left active/selected)
$e3.Range("D6").Select() | Out-Null
$blank.Range("E15").Select() | Out-Null
$blank.Activate() | Out-Null
